$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 69629.414
$ws.Range("I62").Value = 94916.164
$ws.Range("J62").Value = 8941.200000000001
$ws.Range("K62").Value = 94916.164
$ws.Range("L62").Value = 8941.200000000001
$ws.Range("M62").Value = -94292.164
$ws.Range("N62").Value = -10189.2
$ws.Range("H65").Value = 69629.414
$ws.Range("I65").Value = 94916.164
$ws.Range("J65").Value = 8941.200000000001
$ws.Range("K65").Value = 474580.82
$ws.Range("L65").Value = 44706
$ws.Range("M65").Value = -471460.82
$ws.Range("N65").Value = -50946
$ws.Range("H135").Value = 1917.0834
$ws.Range("I135").Value = 833.94446
$ws.Range("J135").Value = 5166.5
$ws.Range("K135").Value = 7505.50014
$ws.Range("L135").Value = 46498.5
$ws.Range("M135").Value = -4970.50014
$ws.Range("N135").Value = -51568.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2352.3489
$ws.Range("I2").Value = 2411.2068
$ws.Range("J2").Value = 2230.4285
$ws.Range("K2").Value = 2411.2068
$ws.Range("L2").Value = 2230.4285
$ws.Range("M2").Value = -2298.2068
$ws.Range("N2").Value = -2456.4285
$ws.Range("H49").Value = 18000
$ws.Range("J49").Value = 18000
$ws.Range("L49").Value = 18000
$ws.Range("N49").Value = -18520
$ws.Range("H53").Value = 4000
$ws.Range("I53").Value = 4000
$ws.Range("K53").Value = 4000
$ws.Range("M53").Value = -3318
$ws.Range("H97").Value = 556.1667
$ws.Range("I97").Value = 536.05554
$ws.Range("J97").Value = 616.5
$ws.Range("K97").Value = 536.05554
$ws.Range("L97").Value = 616.5
$ws.Range("M97").Value = -40.05553999999995
$ws.Range("N97").Value = -1608.5
$ws.Range("H110").Value = 407.3125
$ws.Range("I110").Value = 437.7857
$ws.Range("J110").Value = 194
$ws.Range("K110").Value = 437.7857
$ws.Range("L110").Value = 194
$ws.Range("M110").Value = 1607.2143
$ws.Range("N110").Value = -4284
$ws.Range("H116").Value = 2352.3489
$ws.Range("I116").Value = 2411.2068
$ws.Range("J116").Value = 2230.4285
$ws.Range("K116").Value = 2411.2068
$ws.Range("L116").Value = 2230.4285
$ws.Range("M116").Value = -117.2067999999999
$ws.Range("N116").Value = -6818.4285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2352.3489
$ws.Range("I3").Value = 2411.2068
$ws.Range("J3").Value = 2230.4285
$ws.Range("K3").Value = 2411.2068
$ws.Range("L3").Value = 2230.4285
$ws.Range("M3").Value = -2297.2068
$ws.Range("N3").Value = -2458.4285
$ws.Range("H107").Value = 10194
$ws.Range("I107").Value = 789.63635
$ws.Range("J107").Value = 44676.668
$ws.Range("K107").Value = 789.63635
$ws.Range("L107").Value = 44676.668
$ws.Range("M107").Value = 1130.36365
$ws.Range("N107").Value = -48516.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H131").Value = 26326
$ws.Range("J131").Value = 26326
$ws.Range("L131").Value = 26326
$ws.Range("N131").Value = -36406
$ws.Range("H134").Value = 2266.3333
$ws.Range("I134").Value = 1838.8948
$ws.Range("J134").Value = 3004.6365
$ws.Range("K134").Value = 5516.6844
$ws.Range("L134").Value = 9013.9095
$ws.Range("M134").Value = -2981.6844
$ws.Range("N134").Value = -14083.9095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 2472.7693
$ws.Range("J117").Value = 2655.5
$ws.Range("L117").Value = 7966.5
$ws.Range("N117").Value = -14850.5
$ws.Range("H129").Value = 2393.8572
$ws.Range("I129").Value = 1377.1818
$ws.Range("J129").Value = 3051.7058
$ws.Range("K129").Value = 4131.5454
$ws.Range("L129").Value = 9155.117400000001
$ws.Range("M129").Value = 868.4546
$ws.Range("N129").Value = -19155.1174

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3170
$ws.Range("I80").Value = 3334
$ws.Range("J80").Value = 2350
$ws.Range("K80").Value = 3334
$ws.Range("L80").Value = 2350
$ws.Range("M80").Value = -2336
$ws.Range("N80").Value = -4346
$ws.Range("H83").Value = 3170
$ws.Range("I83").Value = 3334
$ws.Range("J83").Value = 2350
$ws.Range("K83").Value = 16670
$ws.Range("L83").Value = 11750
$ws.Range("M83").Value = -11678
$ws.Range("N83").Value = -21734
$ws.Range("H122").Value = 2315.2188
$ws.Range("I122").Value = 1757.3158
$ws.Range("J122").Value = 3130.6155
$ws.Range("K122").Value = 5271.9474
$ws.Range("L122").Value = 9391.8465
$ws.Range("M122").Value = -2821.9474
$ws.Range("N122").Value = -14291.8465

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1525.0834
$ws.Range("I16").Value = 1711.2222
$ws.Range("J16").Value = 966.6667
$ws.Range("K16").Value = 1711.2222
$ws.Range("L16").Value = 966.6667
$ws.Range("M16").Value = -1541.2222
$ws.Range("N16").Value = -1306.6667
$ws.Range("H41").Value = 13000
$ws.Range("J41").Value = 13000
$ws.Range("L41").Value = 13000
$ws.Range("N41").Value = -13876
$ws.Range("H47").Value = 13475
$ws.Range("J47").Value = 13475
$ws.Range("L47").Value = 13475
$ws.Range("N47").Value = -14455
$ws.Range("H48").Value = 3041
$ws.Range("I48").Value = 3041
$ws.Range("K48").Value = 3041
$ws.Range("M48").Value = -2380
$ws.Range("H52").Value = 13475
$ws.Range("J52").Value = 13475
$ws.Range("L52").Value = 13475
$ws.Range("N52").Value = -13941
$ws.Range("H61").Value = 2662.8462
$ws.Range("I61").Value = 2311.7
$ws.Range("J61").Value = 3833.3333
$ws.Range("K61").Value = 2311.7
$ws.Range("L61").Value = 3833.3333
$ws.Range("M61").Value = -2109.7
$ws.Range("N61").Value = -4237.3333
$ws.Range("H100").Value = 7015711
$ws.Range("I100").Value = 8017655
$ws.Range("K100").Value = 8017655
$ws.Range("M100").Value = -8017114
$ws.Range("H113").Value = 2662.8462
$ws.Range("I113").Value = 2311.7
$ws.Range("J113").Value = 3833.3333
$ws.Range("K113").Value = 2311.7
$ws.Range("L113").Value = 3833.3333
$ws.Range("M113").Value = -141.6999999999998
$ws.Range("N113").Value = -8173.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1408.7567
$ws.Range("I136").Value = 1429.3636
$ws.Range("J136").Value = 1238.75
$ws.Range("K136").Value = 4288.0908
$ws.Range("L136").Value = 3716.25
$ws.Range("M136").Value = -1738.0908
$ws.Range("N136").Value = -8816.25
